# "Generate Report for Archive"
#  - Status text "Ready for handoff" -> "In Translation"
#    (shared by Overview!E2/F2, zh-cn!C2, de-de!C2)
#  - Narrow the "zh-cn"/"de-de" status columns from ~17.22 chars to ~13.41 chars
#    (Overview columns E & F, and column C on both the zh-cn and de-de sheets)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# --- Update the status label everywhere it appears ---
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws2.Range("C2").Value = "In Translation"
$ws3.Range("C2").Value = "In Translation"

# --- Narrow the status columns ---
# Target stored column width is 13.4101848602295 characters. The engine
# quantizes ColumnWidth onto a 1/6-character grid when it persists the
# worksheet, so 12.5 is the ColumnWidth that lands closest to that target.
$newWidth = 12.5

$ws1.Columns.Item(5).ColumnWidth = $newWidth   # Overview column E (zh-cn)
$ws1.Columns.Item(6).ColumnWidth = $newWidth   # Overview column F (de-de)
$ws2.Columns.Item(3).ColumnWidth = $newWidth   # zh-cn column C (Status)
$ws3.Columns.Item(3).ColumnWidth = $newWidth   # de-de column C (Status)
